# Update Name of Algo
# Refresh imputed values for the RandomForest result sheet (columns A-E,
# rows 2-102) to reflect the latest algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.92220000000002
$ws.Range("C3").Value = -11.9569
$ws.Range("E8").Value = 16.476
$ws.Range("D19").Value = -8.669699999999995
$ws.Range("A21").Value = -20.18029999999998
$ws.Range("A23").Value = -20.01409999999997
$ws.Range("E23").Value = 16.32769999999999
$ws.Range("C24").Value = -13.02579999999999
$ws.Range("D24").Value = -7.8547
$ws.Range("A25").Value = -21.7774
$ws.Range("E26").Value = 15.99599999999999
$ws.Range("B27").Value = 6.403800000000006
$ws.Range("D30").Value = -7.415500000000001
$ws.Range("B31").Value = 5.746500000000003
$ws.Range("D31").Value = -8.218599999999997
$ws.Range("D33").Value = -7.992899999999995
$ws.Range("E37").Value = 16.55770000000001
$ws.Range("B39").Value = 10.0873
$ws.Range("B48").Value = 5.596900000000006
$ws.Range("E48").Value = 17.34240000000001
$ws.Range("B51").Value = 5.512299999999999
$ws.Range("B52").Value = 5.088800000000003
$ws.Range("A53").Value = -21.8196
$ws.Range("B55").Value = 5.892799999999997
$ws.Range("D55").Value = -8.110599999999998
$ws.Range("B56").Value = 5.241199999999997
$ws.Range("A57").Value = -21.92759999999999
$ws.Range("B57").Value = 5.440699999999998
$ws.Range("C57").Value = -12.84789999999999
$ws.Range("A59").Value = -22.3641
$ws.Range("C61").Value = -13.01249999999999
$ws.Range("E62").Value = 16.4659
$ws.Range("D65").Value = -8.145299999999997
$ws.Range("E66").Value = 17.07660000000002
$ws.Range("A69").Value = -21.6649
$ws.Range("C70").Value = -12.068
$ws.Range("D70").Value = -8.393899999999997
$ws.Range("B73").Value = 8.614799999999995
$ws.Range("D75").Value = -8.167999999999999
$ws.Range("A79").Value = -20.60810000000001
$ws.Range("A83").Value = -21.9881
$ws.Range("D83").Value = -8.343399999999994
$ws.Range("C86").Value = -13.00079999999999
$ws.Range("B89").Value = 4.993699999999993
$ws.Range("E89").Value = 17.18020000000002
$ws.Range("B90").Value = 5.302299999999999
$ws.Range("A93").Value = -21.25229999999999
$ws.Range("E94").Value = 18.96560000000002
$ws.Range("D96").Value = -7.5106
$ws.Range("D97").Value = -8.411999999999999
$ws.Range("C98").Value = -12.11869999999999
$ws.Range("C100").Value = -12.47919999999999
$ws.Range("C102").Value = -13.4235
